# "add another feature file" - adds a second test-case sheet (TC_02) with
# Selenium-style locator data, following the same Name/Locator/Value layout
# as the existing TC_01 sheet, and makes it the active sheet/tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Sheet1 (TC_01) keeps its previous selection, but since it stops being the
# active tab once TC_02 is added, move the cursor to where the author last
# left it on that sheet (B10) before switching sheets.
$null = $ws1.Range("B10").Select()

# Insert the new sheet right after TC_01 and name it TC_02.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "TC_02"

# Column widths roughly matching the authored sheet (A/C ~46, B ~45).
$ws2.Columns.Item(1).ColumnWidth = 46.12
$ws2.Columns.Item(2).ColumnWidth = 45.3
$ws2.Columns.Item(3).ColumnWidth = 46.12

# Header row - identical to TC_01's header.
$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Locator"
$ws2.Range("C1").Value = "Value"

# Row 2 - search box / search term
$ws2.Range("A2").Value = "txt_search_box"
$ws2.Range("B2").Value = "//input[@id='search_query_top']"
$ws2.Range("C2").Value = "Printed"

# Row 3 - search button
$ws2.Range("A3").Value = "btn_search"
$ws2.Range("B3").Value = "//button[@name='submit_search']"

# Row 4 - product link
$ws2.Range("A4").Value = "btn_product_name"
$ws2.Range("B4").Value = "//img[@title='Printed Chiffon Dress']"

# Row 5 - quantity / plus icon
$ws2.Range("A5").Value = "icon_plus"
$ws2.Range("B5").Value = "//i[@class='icon-plus']"

# Row 6 - size selector
$ws2.Range("A6").Value = "select_size"
$ws2.Range("B6").Value = "//select[@id='group_1']"
$ws2.Range("C6").Value = "M"

# Row 7 - color swatch
$ws2.Range("A7").Value = "btn_color"
$ws2.Range("B7").Value = "//a[@id='color_15']"

# Row 8 - add to cart (reuses the existing shared string "btn_submit")
$ws2.Range("A8").Value = "btn_submit"
$ws2.Range("B8").Value = "//span[contains(text(),'Add to cart')]"

# Row 9 - proceed to checkout
$ws2.Range("A9").Value = "btn_proceed_to_checkout"
$ws2.Range("B9").Value = "//a[@title='Proceed to checkout']"

# Row 10 - login e-mail (kept as plain text so it isn't re-interpreted)
$ws2.Range("A10").Value = "txt_email"
$ws2.Range("B10").Value = "//input[@id='email']"
$ws2.Range("C10").NumberFormat = "@"
$ws2.Range("C10").Value = "18naduni@gmail.com"

# Row 11 - login password
$ws2.Range("A11").Value = "txt_password"
$ws2.Range("B11").Value = "//input[@id='passwd']"
$ws2.Range("C11").Value = "hnUg2kaF@MaUrfr"

# Row 12 - login submit
$ws2.Range("A12").Value = "btn_login"
$ws2.Range("B12").Value = "//button[@id='SubmitLogin']"

# Row 13 - checkout step 1
$ws2.Range("A13").Value = "btn_checkout"
$ws2.Range("B13").Value = "//a[@href='http://automationpractice.com/index.php?controller=order&step=1']"

# Row 14 - checkout step 2 (address)
$ws2.Range("A14").Value = "btn_checkout2"
$ws2.Range("B14").Value = "//button[@name='processAddress']"

# Row 15 - agree to terms checkbox
$ws2.Range("A15").Value = "chk_agree"
$ws2.Range("B15").Value = "//div[@id='uniform-cgv']"

# Row 16 - checkout step 3 (carrier)
$ws2.Range("A16").Value = "btn_checkout3"
$ws2.Range("B16").Value = "//button[@name='processCarrier']"

# Row 17 - shipping / payment method
$ws2.Range("A17").Value = "btn_shipping"
$ws2.Range("B17").Value = "//a[@title='Pay by bank wire']"

# Row 18 - confirm order
$ws2.Range("A18").Value = "btn_confirm"
$ws2.Range("B18").Value = "//span[text()='I confirm my order']"

# Row 19 - order confirmation text
$ws2.Range("A19").Value = "txt_confirm"
$ws2.Range("B19").Value = "//strong[text()='Your order will be sent as soon as we receive payment.']"

# Leave the cursor where the author left it, and make TC_02 the active tab.
$null = $ws2.Range("B19").Select()
$null = $ws2.Activate()
